$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.896.17"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "1.844.91"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.49%  "

$ws.Range("E8").Value = "  +2.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07180"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9226"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07602"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.25%  "

$ws.Range("D13").Value = "1.833.25"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.391"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008641"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "26.935.11"
$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("E21").Value = "  +2.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.022"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.913"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.006"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.871"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08852"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("E31").Value = "  +4.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7470"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.796"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.166"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.486"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.089"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05262"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.20%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01949"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.959"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5213"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.891"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.49%  "

$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.181"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4694"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.70%  "

$ws.Range("E48").Value = "  +2.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06021"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8842"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.26%  "
